# The sheet originally had three columns: BRANCH_CODE (A), BRANCH_NAME (B)
# and PASSWORD (C). The edit removes the BRANCH_CODE column entirely, so
# BRANCH_NAME becomes column A and PASSWORD becomes column B (with all of
# their data shifting left accordingly). The now-unused "BRANCH_CODE" /
# "1-00x" shared strings and the Arial-font style that was only applied to
# the old column A data cells disappear as a natural consequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A:A").EntireColumn.Delete()

# Match the author's resulting selection state (A1:A1048576 instead of
# the original A1:C1048576).
$ws.Range("A1:A1048576").Select()
